$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reset completeness-mandatory counters for several attributes
$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 0
$ws.Range("B16").Value = 0
$ws.Range("B17").Value = 0
$ws.Range("B18").Value = 0
$ws.Range("B19").Value = 0

# Populate metadata-compliance counter/score for NOTE_ETXT (row 4) and NOTES_FTXT (row 5)
$ws.Range("J4").Value = 57
$ws.Range("K4").Value = 0.014984
$ws.Range("J5").Value = 57
$ws.Range("K5").Value = 0.014984
